# "new run with new features" -- refresh the classifier probability-averaging
# evaluation table with results from a new run (new feature set). Row labels
# get reshuffled (MLP 32/16/64/128 interleaved differently, SVC poly /
# RandomForestClassifier swap, GradientBoostingClassifier / LogisticRegressionCV
# / SVC poly shift down one) and every metric column (B:H) gets new numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: SVC rbf (label unchanged)
$ws.Range("A2").Value = "SVC rbf"
$ws.Range("B2").Value = 0.6511
$ws.Range("C2").Value = 0.6511
$ws.Range("D2").Value = 0.9451
$ws.Range("E2").Value = 0.5214
$ws.Range("F2").Value = 0.4913
$ws.Range("G2").Value = 0.7305
$ws.Range("H2").Value = 0.7304

# Row 3: MLP-deep (label unchanged)
$ws.Range("A3").Value = "MLP-deep"
$ws.Range("B3").Value = 0.6538
$ws.Range("C3").Value = 0.6538
$ws.Range("D3").Value = 0.9409
$ws.Range("E3").Value = 0.5317
$ws.Range("F3").Value = 0.4945
$ws.Range("G3").Value = 0.7302
$ws.Range("H3").Value = 0.7251

# Row 4: was MLP 128, now MLP 32
$ws.Range("A4").Value = "MLP 32"
$ws.Range("B4").Value = 0.5989
$ws.Range("C4").Value = 0.5989
$ws.Range("D4").Value = 0.9341
$ws.Range("E4").Value = 0.5824
$ws.Range("F4").Value = 0.5428
$ws.Range("G4").Value = 0.6996
$ws.Range("H4").Value = 0.6988

# Row 5: was MLP 32, now MLP 16
$ws.Range("A5").Value = "MLP 16"
$ws.Range("B5").Value = 0.5755
$ws.Range("C5").Value = 0.5755
$ws.Range("D5").Value = 0.9327
$ws.Range("E5").Value = 0.5926
$ws.Range("F5").Value = 0.549
$ws.Range("G5").Value = 0.6935
$ws.Range("H5").Value = 0.6935

# Row 6: MLP 64 (label unchanged)
$ws.Range("A6").Value = "MLP 64"
$ws.Range("B6").Value = 0.6003
$ws.Range("C6").Value = 0.6003
$ws.Range("D6").Value = 0.9327
$ws.Range("E6").Value = 0.5933
$ws.Range("F6").Value = 0.5413
$ws.Range("G6").Value = 0.6936
$ws.Range("H6").Value = 0.6932

# Row 7: was GradientBoostingClassifier, now MLP 128
$ws.Range("A7").Value = "MLP 128"
$ws.Range("B7").Value = 0.5618
$ws.Range("C7").Value = 0.5618
$ws.Range("D7").Value = 0.9272
$ws.Range("E7").Value = 0.6113
$ws.Range("F7").Value = 0.5628
$ws.Range("G7").Value = 0.6864
$ws.Range("H7").Value = 0.6839

# Row 8: was LogisticRegressionCV, now GradientBoostingClassifier
$ws.Range("A8").Value = "GradientBoostingClassifier"
$ws.Range("B8").Value = 0.581
$ws.Range("C8").Value = 0.581
$ws.Range("D8").Value = 0.9258
$ws.Range("E8").Value = 0.6101
$ws.Range("F8").Value = 0.5532
$ws.Range("G8").Value = 0.6846
$ws.Range("H8").Value = 0.6845

# Row 9: was MLP 16, now LogisticRegressionCV
$ws.Range("A9").Value = "LogisticRegressionCV"
$ws.Range("B9").Value = 0.5453
$ws.Range("C9").Value = 0.5453
$ws.Range("D9").Value = 0.9231
$ws.Range("E9").Value = 0.6372
$ws.Range("F9").Value = 0.5926
$ws.Range("G9").Value = 0.6709
$ws.Range("H9").Value = 0.6705

# Row 10: was RandomForestClassifier, now SVC poly
$ws.Range("A10").Value = "SVC poly"
$ws.Range("B10").Value = 0.5041
$ws.Range("C10").Value = 0.5041
$ws.Range("D10").Value = 0.9052
$ws.Range("E10").Value = 0.7045
$ws.Range("F10").Value = 0.6256
$ws.Range("G10").Value = 0.6356
$ws.Range("H10").Value = 0.6356

# Row 11: was SVC poly, now RandomForestClassifier
$ws.Range("A11").Value = "RandomForestClassifier"
$ws.Range("B11").Value = 0.4286
$ws.Range("C11").Value = 0.4286
$ws.Range("D11").Value = 0.8997
$ws.Range("E11").Value = 0.8391
$ws.Range("F11").Value = 0.7178
$ws.Range("G11").Value = 0.5661
$ws.Range("H11").Value = 0.5661

# Row 12: SVC sigmoid (label unchanged)
$ws.Range("A12").Value = "SVC sigmoid"
$ws.Range("B12").Value = 0.4739
$ws.Range("C12").Value = 0.4739
$ws.Range("D12").Value = 0.8791
$ws.Range("E12").Value = 0.9139
$ws.Range("F12").Value = 0.7194
$ws.Range("G12").Value = 0.5279
$ws.Range("H12").Value = 0.5274
